$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("identifiers_extended")

# Update rows 55-56 (existing rows get new A value + new B value)
# Add new rows 57-115 with Code/Chemogenomic_label pairs
$ws.Cells.Item(55, 1).Value = 'CRL'
$ws.Cells.Item(55, 2).Value = 'CERULENIN-6.0'
$ws.Cells.Item(56, 1).Value = 'C90'
$ws.Cells.Item(56, 2).Value = 'CHIR090-0.075'
$ws.Cells.Item(57, 1).Value = 'CPZ'
$ws.Cells.Item(57, 2).Value = 'CHLOROPROMAZINE-24'
$ws.Cells.Item(58, 1).Value = 'CHO'
$ws.Cells.Item(58, 2).Value = 'CHOLATE-2.0%'
$ws.Cells.Item(59, 1).Value = 'CPT'
$ws.Cells.Item(59, 2).Value = 'CISPLATIN-100'
$ws.Cells.Item(60, 1).Value = 'CSD'
$ws.Cells.Item(60, 2).Value = 'CYCLOSERINED-16'
$ws.Cells.Item(61, 1).Value = 'DCH'
$ws.Cells.Item(61, 2).Value = 'DEOXYCHOLATE-2.0%'
$ws.Cells.Item(62, 1).Value = 'DBC'
$ws.Cells.Item(62, 2).Value = 'DIBUCAINE-1.2'
$ws.Cells.Item(63, 1).Value = 'DXR'
$ws.Cells.Item(63, 2).Value = 'DOXORUBICIN-10.0'
$ws.Cells.Item(64, 1).Value = 'EDTA'
$ws.Cells.Item(64, 2).Value = 'EDTA-1.0'
$ws.Cells.Item(65, 1).Value = 'EGCG'
$ws.Cells.Item(65, 2).Value = 'EGCG-50'
$ws.Cells.Item(66, 1).Value = 'EGTA'
$ws.Cells.Item(66, 2).Value = 'EGTA-2.0'
$ws.Cells.Item(67, 1).Value = 'EPI'
$ws.Cells.Item(67, 2).Value = 'EPINEPHRINE-1000'
$ws.Cells.Item(68, 1).Value = 'ETH'
$ws.Cells.Item(68, 2).Value = 'ETHANOL-6.0'
$ws.Cells.Item(69, 1).Value = 'EDB'
$ws.Cells.Item(69, 2).Value = 'ETHIDIUMBROMIDE-50'
$ws.Cells.Item(70, 1).Value = 'GLUCOSAMINE'
$ws.Cells.Item(70, 2).Value = 'GLUCOSAMINE'
$ws.Cells.Item(71, 1).Value = 'GLUCOSE'
$ws.Cells.Item(71, 2).Value = 'GLUCOSE'
$ws.Cells.Item(72, 1).Value = 'GFOS'
$ws.Cells.Item(72, 2).Value = 'GLUFOSFOMYCIN-0.2'
$ws.Cells.Item(73, 1).Value = 'GLYCEROL'
$ws.Cells.Item(73, 2).Value = 'GLYCEROL'
$ws.Cells.Item(74, 1).Value = 'HCO'
$ws.Cells.Item(74, 2).Value = 'HIGHCOBALT-0.5'
$ws.Cells.Item(75, 1).Value = 'HCU'
$ws.Cells.Item(75, 2).Value = 'HIGHCOPPER-4.0'
$ws.Cells.Item(76, 1).Value = 'HFE'
$ws.Cells.Item(76, 2).Value = 'HIGHFE'
$ws.Cells.Item(77, 1).Value = 'HNI'
$ws.Cells.Item(77, 2).Value = 'HIGHNICKEL-1.0'
$ws.Cells.Item(78, 1).Value = 'HUREA'
$ws.Cells.Item(78, 2).Value = 'HYDROXYUREA-10.0'
$ws.Cells.Item(79, 1).Value = 'INDO'
$ws.Cells.Item(79, 2).Value = 'INDOLICIDIN-0.1'
$ws.Cells.Item(80, 1).Value = 'INZ'
$ws.Cells.Item(80, 2).Value = 'ISONIAZID-1.5'
$ws.Cells.Item(81, 1).Value = 'LFE'
$ws.Cells.Item(81, 2).Value = 'LOWFE'
$ws.Cells.Item(82, 1).Value = 'MALTOSE'
$ws.Cells.Item(82, 2).Value = 'MALTOSE'
$ws.Cells.Item(83, 1).Value = 'MEC'
$ws.Cells.Item(83, 2).Value = 'MECILLINAM-0.12'
$ws.Cells.Item(84, 1).Value = 'MTX'
$ws.Cells.Item(84, 2).Value = 'METHOTREXATE-25'
$ws.Cells.Item(85, 1).Value = 'MMC'
$ws.Cells.Item(85, 2).Value = 'MITOMYCINC-0.1'
$ws.Cells.Item(86, 1).Value = 'MMS'
$ws.Cells.Item(86, 2).Value = 'MMS-0.05%'
$ws.Cells.Item(87, 1).Value = 'NAG'
$ws.Cells.Item(87, 2).Value = 'N-ACETYLGLUCOSAMINE'
$ws.Cells.Item(88, 1).Value = 'NACL'
$ws.Cells.Item(88, 2).Value = 'NACL-600'
$ws.Cells.Item(89, 1).Value = 'NH4CL'
$ws.Cells.Item(89, 2).Value = 'NH4CL'
$ws.Cells.Item(90, 1).Value = 'NIG'
$ws.Cells.Item(90, 2).Value = 'NIGERICIN-5.0'
$ws.Cells.Item(91, 1).Value = 'NEPI'
$ws.Cells.Item(91, 2).Value = 'NOREPINEPHRINE-1000'
$ws.Cells.Item(92, 1).Value = 'NOV'
$ws.Cells.Item(92, 2).Value = 'NOVOBIOCIN-30'
$ws.Cells.Item(93, 1).Value = 'PQ'
$ws.Cells.Item(93, 2).Value = 'PARAQUAT-18.0'
$ws.Cells.Item(94, 1).Value = 'PH4'
$ws.Cells.Item(94, 2).Value = 'PH4'
$ws.Cells.Item(95, 1).Value = 'PH10'
$ws.Cells.Item(95, 2).Value = 'PH10'
$ws.Cells.Item(96, 1).Value = 'PHL'
$ws.Cells.Item(96, 2).Value = 'PHLEOMYCIN-1.0'
$ws.Cells.Item(97, 1).Value = 'PMS'
$ws.Cells.Item(97, 2).Value = 'PMS-0.1'
$ws.Cells.Item(98, 1).Value = 'PRO'
$ws.Cells.Item(98, 2).Value = 'PROCAINE-30'
$ws.Cells.Item(99, 1).Value = 'PPI'
$ws.Cells.Item(99, 2).Value = 'PROPIDIUMIODIDE-50'
$ws.Cells.Item(100, 1).Value = 'PUR'
$ws.Cells.Item(100, 2).Value = 'PUROMYCIN-25'
$ws.Cells.Item(101, 1).Value = 'PYO'
$ws.Cells.Item(101, 2).Value = 'PYOCYANIN-10.0'
$ws.Cells.Item(102, 1).Value = 'RAD'
$ws.Cells.Item(102, 2).Value = 'RADICICOL-10'
$ws.Cells.Item(103, 1).Value = 'SDS'
$ws.Cells.Item(103, 2).Value = 'SDS-4.0%'
$ws.Cells.Item(104, 1).Value = 'SDSEDTA'
$ws.Cells.Item(104, 2).Value = 'SDS1.0%/EDTA0.5'
$ws.Cells.Item(105, 1).Value = 'STN'
$ws.Cells.Item(105, 2).Value = 'SDS1.0%/EDTA0.5'
$ws.Cells.Item(106, 1).Value = 'SUCCINATE'
$ws.Cells.Item(106, 2).Value = 'SUCCINATE'
$ws.Cells.Item(107, 1).Value = 'SMZ'
$ws.Cells.Item(107, 2).Value = 'SULFAMETHIZOLE-300'
$ws.Cells.Item(108, 1).Value = 'TCHO'
$ws.Cells.Item(108, 2).Value = 'TAUROCHOLATE-1.0%'
$ws.Cells.Item(109, 1).Value = 'THP'
$ws.Cells.Item(109, 2).Value = 'THEOPHYLLINE-100'
$ws.Cells.Item(110, 1).Value = 'TLM'
$ws.Cells.Item(110, 2).Value = 'THIOLACTOMYCIN-50'
$ws.Cells.Item(111, 1).Value = 'TMPSMZ'
$ws.Cells.Item(111, 2).Value = 'TRIMETHOPRIM-0.1,SULFAMETHIZOLE-50'
$ws.Cells.Item(112, 1).Value = 'TTX'
$ws.Cells.Item(112, 2).Value = 'TRITONX-0.2%'
$ws.Cells.Item(113, 1).Value = 'TUN'
$ws.Cells.Item(113, 2).Value = 'TUNICAMYCIN-7.5'
$ws.Cells.Item(114, 1).Value = 'UV'
$ws.Cells.Item(114, 2).Value = 'UV-24SEC'
$ws.Cells.Item(115, 1).Value = 'VERA'
$ws.Cells.Item(115, 2).Value = 'VERAPAMIL-1.0'

# Update the view: scroll/select to match final state
$ws.Activate() | Out-Null
$ws.Range("B115").Select() | Out-Null
